$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column B ("Ma dat hang" / order code column).
# This shifts every existing column from B..V to C..W and Excel auto-adjusts
# all formulas/merged ranges/column widths along the way.
$ws.Columns.Item(2).Insert()

# Header text for the newly inserted column.
$ws.Range("B4").Value = "Mã đặt hàng"

# Row 29 (totals row): the former I/J "Tong tien truoc CK" / "Tong chiet khau"
# SUM cells (now shifted to J/K) get a right-aligned look instead of the
# default centered one used by the rest of the totals row.
$ws.Range("J29:K29").HorizontalAlignment = -4152  ' xlRight

# B29 sits inside the continuous outlined footer box together with A29; after
# the column insert it should not carry a left border (that border belongs to
# A29 only), matching the rest of the interior footer cells.
$ws.Range("B29").Borders.Item(7).LineStyle = -4142  ' xlLineStyleNone (xlEdgeLeft = 7)

# The "Tong chi nhanh" / branch-name data column (old I, now J) and the
# "Tong tien truoc CK" data column (old J, now K) switch from center to
# right aligned text for their data rows.
$ws.Range("J5:K28").HorizontalAlignment = -4152  ' xlRight
